# Apply updated values ("Results from R script") to the
# "Crédito disponível - Centralização - Campus Itabaiana" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("N11").Value = 14758.31
$ws.Range("O11").Value = 14758.31

# Row 13
$ws.Range("N13").Value = 608625.34

# Row 14
$ws.Range("O14").Value = 221318.08

# Row 15
$ws.Range("O15").Value = 26968.24

# Row 20
$ws.Range("N20").Value = 39147.83

# Row 28
$ws.Range("K28").Value = 76128.56

# Row 30 (previously empty cell gets a value)
$ws.Range("N30").Value = 2991.3
